# Plantilla Lista de Tareas de la 5ta Iteración
# Registrar el consumo de horas del Día 1 (columna H) para la primera
# tarea de la iteración (fila 6). El resto de columnas de la fila (horas
# restantes, acumulados, etc.) están calculadas con fórmulas que se
# recalculan automáticamente a partir de este valor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Horas consumidas el Día 1 en la tarea de la fila 6
$ws.Range("H6").Value = 2

# Deja la selección / vista tal como quedó al guardar el archivo
$ws.Range("K8").Select()
